$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Merge the two runs "Umstrukturieren Befragungserstellung" and
#    " (Krukenfellner)" into a single run with the combined text.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Umstrukturieren Befragungserstellung (Krukenfellner)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Umstrukturieren Befragungserstellung (Krukenfellner)", 2)

# ---------------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark from the empty paragraph after
#    "Fortsetzung Gruppenverwaltung (Krukenfellner)" to the very start of the
#    document (right before the "Besprechungsprotokoll" heading run).
# ---------------------------------------------------------------------------

# Remove the bookmark from its current location.
$oldBm = $d.Bookmarks("_GoBack")
$oldBm.Delete()

# Re-create it at the very beginning of the document. A zero-length range
# placed exactly at position 0 tends to anchor to the following paragraph,
# so a throw-away marker character is used to get a reliable anchor, and is
# removed again afterwards.
$startRange = $d.Range(0, 0)
$startRange.InsertBefore("X")
$markerRange = $d.Range(1, 1)
$d.Bookmarks.Add("_GoBack", $markerRange)
$markerChar = $d.Range(0, 1)
$markerChar.Text = ""

# ---------------------------------------------------------------------------
# 3) Update the cached result of the header's TIME field from 2018-11-07 to
#    2018-11-19.
# ---------------------------------------------------------------------------
$header = $d.Sections.Item(1).Headers.Item(1)
$header.Range.Find.Execute(
    "2018-11-07", $true, $false, $false, $false, $false, $true, 1, $false,
    "2018-11-19", 2)
